# Add "Organ Confined Disease" (OCD) column to the dataset, between SVI (N) and
# METASTASIS (old O). This mirrors: select column O, insert a new blank column
# (shifting METASTASIS..DEATH_TIME one column to the right), then fill in the
# new OCD header + values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new, empty column at O - this shifts the existing O:T data
#    (METASTASIS, METASTASIS_TIME, CRPC, CRPC_TIME, DEATH, DEATH_TIME) right to P:U,
#    carrying their values/formatting with them.
$ws.Range("O1").EntireColumn.Insert()

# 2) Header for the newly inserted column.
$ws.Range("O1").Value = "OCD"

# 3) Organ-confined-disease flag per patient row (2-7).
$ws.Range("O2").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("O7").Value = 1

# 4) A couple of data corrections made alongside the column insert.
$ws.Range("U4").Value = 140
$ws.Range("T7").Value = 1

# 5) Column width tweaks that accompanied the edit.
$ws.Range("K1").ColumnWidth = 10.140625
$ws.Range("O1").ColumnWidth = 6
$ws.Range("R1").ColumnWidth = 5.7109375

# 6) Leave the selection where the author last left it.
$ws.Range("T13").Select()
